# The sheet gained one more weekly data row. A new record (most recent
# week) is inserted at row 52, pushing the previously-existing rows 52-117
# down to 53-118 (dimension grows from A1:R117 to A1:R118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; everything from the old row 52 onward
# shifts down by one row automatically.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(52, 1).Value2  = 6
$ws.Cells.Item(52, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(52, 3).Value2  = "Metropolitana"
$ws.Cells.Item(52, 4).Value2  = 44495
$ws.Cells.Item(52, 5).Value2  = 13
$ws.Cells.Item(52, 6).Value2  = 100112029
$ws.Cells.Item(52, 7).Value2  = "Orégano"
$ws.Cells.Item(52, 8).Value2  = "Sin especificar"
$ws.Cells.Item(52, 9).Value2  = "Primera"
$ws.Cells.Item(52, 10).Value2 = 35
$ws.Cells.Item(52, 11).Value2 = 8000
$ws.Cells.Item(52, 12).Value2 = 9000
$ws.Cells.Item(52, 13).Value2 = 8457
$ws.Cells.Item(52, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(52, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value2 = 2819
$ws.Cells.Item(52, 17).Value2 = 3
$ws.Cells.Item(52, 18).Value2 = "Hortaliza"
